$d = $word.ActiveDocument

# --- 1. 'Wish List' -> 'Brain Dump', with ListParagraph style applied ---
$found = $d.Content.Find.Execute("Wish List", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "could not find 'Wish List' paragraph" }
$wishPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Wish List") {
        $wishPara = $p
        break
    }
}
if ($null -eq $wishPara) { throw "could not locate 'Wish List' paragraph object" }
$wishRange = $wishPara.Range
$wishRange.End = $wishRange.End - 1
$wishRange.Text = "Brain Dump"
$wishPara.Style = "List Paragraph"

# --- 2. Append 15 new bullet paragraphs after the '... only' list item ---
$lastListPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -like "*buffs/debuffs only") {
        $lastListPara = $p
    }
}
if ($null -eq $lastListPara) { throw "could not locate the '... only' paragraph" }
# Insert one character before the paragraph mark (i.e. before the trailing
# glyph of ' only') rather than exactly at the paragraph boundary: InsertXML
# at an exact paragraph boundary merges its *last* fragment paragraph into the
# following paragraph (eating the blank w:p that must stay untouched at the
# end of the document); splitting one character earlier avoids that merge and
# Word rejoins the identically-formatted run halves back together.
$splitPoint = $lastListPara.Range.End - 1
$insertAt = $d.Range($splitPoint, $splitPoint)

$newParasXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Spell availability – equal element, sufficient rank. Full ranks transfer during class changes</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Class changes – </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>what’s</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> the penalty?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Status effects spread between adjacent units -&gt; kamikaze units?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>AOE skills are impassible, making enlargement skills useful. “whirlwind” skill might stop a character trying to penetrate to the ranged ranks. Spell casters create bigger target while casting (multiply enemy accuracy)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Items of interest show up during quest battles. Going for the item means that character </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>isn’t</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> battling</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Could have waves of enemies while moving through the dungeon, encountering items and other tasks along the way.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Spell damage = wit roll each turn, Damage over time – insight roll at cast</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Hit with spell of same element -&gt; heal</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>“Damage” depends on element, not all “damage” is HP</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Lazy alchemy -&gt; all combos are run through at dungeon exit</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Alchemy loot box -&gt; components determine chances of getting one thing or another, but random</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Culture motivates non-battle actions. Some resources give more XP than others to a character of a given culture. Impact is determined by need of character’s culture’s needs. Changes with global economy.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Can do alchemy during quest run, gaining Item but losing potential XP</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Characters collect pay each week. If you </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>can’t</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> pay, they go to a neighboring territory, becoming an enemy.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Build a gambling system with odds based on weeks played, average level, etc.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertAt.InsertXML($newParasXml)

Write-Host "edit complete"
